$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The time-report table currently ends with a totals ("Yht") row at row 26.
# We need to:
#   1. fix the date in row 25 (it was one day off),
#   2. push the totals row down to row 30,
#   3. fill rows 26-29 with four new days of work (date / hours / description),
#   4. update the SUM formula range to include the new rows.
# ---------------------------------------------------------------------------

# 1) Row 25's date was 45347 (2024-02-25) but should be 45346 (2024-02-24).
$ws.Range("B25").Value = 45346

# 2) Move the totals row (currently B26:D26) down to B30:D30, carrying its
#    formatting (style s=2) along with it. We'll fix the formula afterwards.
$ws.Range("B26:D26").Copy($ws.Range("B30:D30"))

# 3) Seed rows 26-29 with the same cell formatting pattern used by the rows
#    above them (date style s=5, hours style s=3, description style s=4),
#    then overwrite the actual values/text.
$ws.Range("B25:D25").Copy($ws.Range("B26:D26"))
$ws.Range("B25:D25").Copy($ws.Range("B27:D27"))
$ws.Range("B25:D25").Copy($ws.Range("B28:D28"))
$ws.Range("B25:D25").Copy($ws.Range("B29:D29"))

# Row 26: 2024-02-25
$ws.Range("B26").Value = 45347
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = "Katsoin systeminformation työakalua taisin löytää tavan implementoida sen pitää ottaa selvää paremmin"
$ws.Rows(26).RowHeight = 37.5

# Row 27: 2024-02-26
$ws.Range("B27").Value = 45348
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = "Sain systeminformation toimimaan perus tasolla mutta sen pöivityminen ja tyylitely on vielä työn alla"
$ws.Rows(27).RowHeight = 37.5

# Row 28: 2024-02-27
$ws.Range("B28").Value = 45349
$ws.Range("C28").Value = 8
$ws.Range("D28").Value = "Koitin saada systeminformaation päivitämään oikealla tavalla. Sain jotain edistyksiä mutta paljon vaikeuksia saada päivitäminen pyörimään erilisellä thread:lla"
$ws.Rows(28).RowHeight = 56.25

# Row 29: 2024-02-28
$ws.Range("B29").Value = 45350
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = "Kirjoitin dokumentaatiota. Korjasin suurimalta osalta systeminformation vaikkakin vaatii hiomista."
$ws.Rows(29).RowHeight = 37.5

# 4) Fix up the totals row now sitting at row 30: keep the "Yht" label, widen
#    the SUM range to cover the newly-added rows, and restore its height.
$ws.Range("C30").Formula = "=SUM(C6:C29)"
$ws.Rows(30).RowHeight = 18.75

# ---------------------------------------------------------------------------
# View state: the sheet was scrolled down a couple of rows and the selection
# moved from I29 to E29.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("E29").Select()
